{"js": "const pairs = [\n  [\"2024-12-12 Thursday\", \"2024-12-13 Friday\"],\n  [\"377\u00d79=3393\", \"894\u00d78=7152\"],\n  [\"547\u00d74=2188\", \"422\u00d72=844\"],\n  [\"391\u00d78=3128\", \"671\u00d72=1342\"],\n  [\"517\u00d73=1551\", \"844\u00d72=1688\"],\n  [\"794\u00d77=5558\", \"538\u00d75=2690\"],\n  [\"835\u00d74=3340\", \"545\u00d77=3815\"],\n  [\"609\u00d73=1827\", \"231\u00d79=2079\"],\n  [\"428\u00d73=1284\", \"657\u00d75=3285\"],\n  [\"981\u00d75=4905\", \"561\u00d72=1122\"],\n  [\"137\u00d76=822\", \"802\u00d74=3208\"],\n  [\"780\u00d74=3120\", \"393\u00d77=2751\"],\n  [\"401\u00d75=2005\", \"626\u00d77=4382\"],\n  [\"420\u00d78=3360\", \"861\u00d73=2583\"],\n  [\"371\u00d72=742\", \"876\u00d74=3504\"],\n  [\"184\u00d72=368\", \"849\u00d74=3396\"],\n  [\"359\u00d73=1077\", \"962\u00d79=8658\"],\n  [\"687\u00d72=1374\", \"250\u00d75=1250\"],\n  [\"647\u00d78=5176\", \"448\u00d73=1344\"],\n  [\"143\u00d77=1001\", \"739\u00d73=2217\"],\n  [\"689\u00d76=4134\", \"815\u00d76=4890\"],\n  [\"505\u00d76=3030\", \"143\u00d73=429\"],\n  [\"550\u00d74=2200\", \"297\u00d76=1782\"],\n  [\"872\u00d77=6104\", \"167\u00d74=668\"],\n  [\"704\u00d78=5632\", \"215\u00d72=430\"],\n  [\"797\u00d73=2391\", \"434\u00d72=868\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2024-12-12 Thursday', '2024-12-13 Friday'),\n    @('377\u00d79=3393', '894\u00d78=7152'),\n    @('547\u00d74=2188', '422\u00d72=844'),\n    @('391\u00d78=3128', '671\u00d72=1342'),\n    @('517\u00d73=1551', '844\u00d72=1688'),\n    @('794\u00d77=5558', '538\u00d75=2690'),\n    @('835\u00d74=3340', '545\u00d77=3815'),\n    @('609\u00d73=1827', '231\u00d79=2079'),\n    @('428\u00d73=1284', '657\u00d75=3285'),\n    @('981\u00d75=4905', '561\u00d72=1122'),\n    @('137\u00d76=822', '802\u00d74=3208'),\n    @('780\u00d74=3120', '393\u00d77=2751'),\n    @('401\u00d75=2005', '626\u00d77=4382'),\n    @('420\u00d78=3360', '861\u00d73=2583'),\n    @('371\u00d72=742', '876\u00d74=3504'),\n    @('184\u00d72=368', '849\u00d74=3396'),\n    @('359\u00d73=1077', '962\u00d79=8658'),\n    @('687\u00d72=1374', '250\u00d75=1250'),\n    @('647\u00d78=5176', '448\u00d73=1344'),\n    @('143\u00d77=1001', '739\u00d73=2217'),\n    @('689\u00d76=4134', '815\u00d76=4890'),\n    @('505\u00d76=3030', '143\u00d73=429'),\n    @('550\u00d74=2200', '297\u00d76=1782'),\n    @('872\u00d77=6104', '167\u00d74=668'),\n    @('704\u00d78=5632', '215\u00d72=430'),\n    @('797\u00d73=2391', '434\u00d72=868'),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$find.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}\n"}
